$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.164.37", "  -0.33%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.853.28", "  -1.63%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.013", "  +0.58%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "313.31", "  -0.22%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.009", "  +0.39%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5080", "  -1.20%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3903", "  -0.10%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.08249", "  -1.52%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.107", "  -1.13%  ")
    ,@("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "41.62", "  +0.05%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.186", "  -0.86%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.852.16", "  -1.56%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.19", "  -2.73%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.172", "  -1.67%  ")
    ,@("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.011", "  +0.41%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001096", "  -1.13%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "90.91", "  -0.59%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06673", "  +0.12%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.009", "  +0.39%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "17.47", "  -2.07%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.905", "  -2.83%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.164.04", "  -0.45%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.01", "  -1.53%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.234", "  -1.64%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.063.39", "  -1.47%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "159.41", "  +0.33%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.58", "  -0.46%  ")
    ,@("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.368", "  -6.01%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "125.85", "  +0.20%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1043", "  -2.22%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.022", "  -2.30%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.761", "  -2.24%  ")
    ,@("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.612", "  +0.42%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02408", "  -1.83%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06414", "  -2.56%  ")
    ,@("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "9.029", "  -7.70%  ")
    ,@("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2160", "  -1.57%  ")
    ,@("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.241", "  +0.52%  ")
    ,@("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6373", "  -2.68%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.173", "  -3.36%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.917", "  -2.38%  ")
    ,@("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "11.04", "  -2.50%  ")
    ,@("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5969", "  -2.84%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.680", "  +0.02%  ")
    ,@("WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.273", "  -1.23%  ")
    ,@("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "12.72", "  -3.56%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.961", "  -2.79%  ")
    ,@("EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.197", "  -3.25%  ")
    ,@("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "120.23", "  -1.25%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06810", "  -1.39%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
